# Reorder compound rows: phenol moves to row 2, dodecane shifts to row 3,
# naphthalene shifts to row 4 (values only shift; row 2's style stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("phenol",      "phenol",      "C6H6O",  "C1=CC=C(C=C1)O",       94.11, 1.5, 6,  6,  1, 0.765763468281798, 0.06426522154925088, 0.1700031877590054, 0,  6,  1, 0,                   0.8193178195728402, 0.1807140580172139),
    @("dodecane",    "dodecane",    "C12H26", "CCCCCCCCCCCC",         170.33,6.1, 12, 26, 0, 0.846192684788352, 0.1538660247754359,  0,                   12, 0,  0, 1.000058709563788,  0,                  0),
    @("naphthalene", "naphthalene", "C10H8",  "C1=CC=C2C=CC=CC2=C1",  128.17,3.3, 10, 8,  0, 0.9371147694468284,0.06291643910431459, 0,                   0,  10, 0, 0,                   1.000031208551143, 0)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
